$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1=14, Q1=15 with the same style as the
# existing header cells (copy format from O1, then set the values). ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: swap I<->K and M<->O values, and append new
# columns P and Q (both with value 2, no special style). ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I{r}
    $ws.Cells.Item($r, 11).Value = 1   # K{r}
    $ws.Cells.Item($r, 13).Value = 2   # M{r}
    $ws.Cells.Item($r, 15).Value = 1   # O{r}
    $ws.Cells.Item($r, 16).Value = 2   # P{r}
    $ws.Cells.Item($r, 17).Value = 2   # Q{r}
}
